$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "D" column (截止一直未充电时间) for rows 2-17 to the new refresh timestamp
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 4).Value = 45980.313275462962
}

# Update rows 18-52 (A: site name, B: terminal name, C: last-charge-end time, D: cutoff time)
$ws.Cells.Item(18, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18, 2).Value = "702号直流"
$ws.Cells.Item(18, 3).Value = 45978.583449074074
$ws.Cells.Item(18, 4).Value = 45980.313275462962
$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "503号直流"
$ws.Cells.Item(19, 3).Value = 45978.603483796294
$ws.Cells.Item(19, 4).Value = 45980.313275462962
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "B01号直流"
$ws.Cells.Item(20, 3).Value = 45978.664583333331
$ws.Cells.Item(20, 4).Value = 45980.313275462962
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "905号直流"
$ws.Cells.Item(21, 3).Value = 45979.055289351854
$ws.Cells.Item(21, 4).Value = 45980.313275462962
$ws.Cells.Item(22, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22, 2).Value = "402号直流"
$ws.Cells.Item(22, 3).Value = 45979.18167824074
$ws.Cells.Item(22, 4).Value = 45980.313275462962
$ws.Cells.Item(23, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(23, 2).Value = "101号直流"
$ws.Cells.Item(23, 3).Value = 45979.18986111111
$ws.Cells.Item(23, 4).Value = 45980.313275462962
$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24, 2).Value = "401号直流"
$ws.Cells.Item(24, 3).Value = 45979.245381944442
$ws.Cells.Item(24, 4).Value = 45980.313275462962
$ws.Cells.Item(25, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(25, 2).Value = "A05号直流"
$ws.Cells.Item(25, 3).Value = 45979.384791666664
$ws.Cells.Item(25, 4).Value = 45980.313275462962
$ws.Cells.Item(26, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(26, 2).Value = "204号直流"
$ws.Cells.Item(26, 3).Value = 45979.52306712963
$ws.Cells.Item(26, 4).Value = 45980.313275462962
$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(27, 2).Value = "904号直流"
$ws.Cells.Item(27, 3).Value = 45979.527511574073
$ws.Cells.Item(27, 4).Value = 45980.313275462962
$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(28, 2).Value = "003B号直流"
$ws.Cells.Item(28, 3).Value = 45979.545555555553
$ws.Cells.Item(28, 4).Value = 45980.313275462962
$ws.Cells.Item(29, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(29, 2).Value = "004A号直流"
$ws.Cells.Item(29, 3).Value = 45979.550092592595
$ws.Cells.Item(29, 4).Value = 45980.313275462962
$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(30, 2).Value = "103号直流"
$ws.Cells.Item(30, 3).Value = 45979.551099537035
$ws.Cells.Item(30, 4).Value = 45980.313275462962
$ws.Cells.Item(31, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value = "703号直流"
$ws.Cells.Item(31, 3).Value = 45979.55667824074
$ws.Cells.Item(31, 4).Value = 45980.313275462962
$ws.Cells.Item(32, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(32, 2).Value = "104号直流"
$ws.Cells.Item(32, 3).Value = 45979.560810185183
$ws.Cells.Item(32, 4).Value = 45980.313275462962
$ws.Cells.Item(33, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(33, 2).Value = "901号直流"
$ws.Cells.Item(33, 3).Value = 45979.563657407409
$ws.Cells.Item(33, 4).Value = 45980.313275462962
$ws.Cells.Item(34, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(34, 2).Value = "A01号直流"
$ws.Cells.Item(34, 3).Value = 45979.565532407411
$ws.Cells.Item(34, 4).Value = 45980.313275462962
$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(35, 2).Value = "402号直流"
$ws.Cells.Item(35, 3).Value = 45979.575370370374
$ws.Cells.Item(35, 4).Value = 45980.313275462962
$ws.Cells.Item(36, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(36, 2).Value = "405号直流"
$ws.Cells.Item(36, 3).Value = 45979.585902777777
$ws.Cells.Item(36, 4).Value = 45980.313275462962
$ws.Cells.Item(37, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(37, 2).Value = "108号直流"
$ws.Cells.Item(37, 3).Value = 45979.586111111108
$ws.Cells.Item(37, 4).Value = 45980.313275462962
$ws.Cells.Item(38, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(38, 2).Value = "210号直流"
$ws.Cells.Item(38, 3).Value = 45979.609259259261
$ws.Cells.Item(38, 4).Value = 45980.313275462962
$ws.Cells.Item(39, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(39, 2).Value = "201号直流"
$ws.Cells.Item(39, 3).Value = 45979.628310185188
$ws.Cells.Item(39, 4).Value = 45980.313275462962
$ws.Cells.Item(40, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(40, 2).Value = "406号直流"
$ws.Cells.Item(40, 3).Value = 45979.666643518518
$ws.Cells.Item(40, 4).Value = 45980.313275462962
$ws.Cells.Item(41, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(41, 2).Value = "209号直流"
$ws.Cells.Item(41, 3).Value = 45979.667731481481
$ws.Cells.Item(41, 4).Value = 45980.313275462962
$ws.Cells.Item(42, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(42, 2).Value = "203号直流"
$ws.Cells.Item(42, 3).Value = 45979.66951388889
$ws.Cells.Item(42, 4).Value = 45980.313275462962
$ws.Cells.Item(43, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(43, 2).Value = "206号直流"
$ws.Cells.Item(43, 3).Value = 45979.673009259262
$ws.Cells.Item(43, 4).Value = 45980.313275462962
$ws.Cells.Item(44, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(44, 2).Value = "107号直流"
$ws.Cells.Item(44, 3).Value = 45979.699432870373
$ws.Cells.Item(44, 4).Value = 45980.313275462962
$ws.Cells.Item(45, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(45, 2).Value = "502号直流"
$ws.Cells.Item(45, 3).Value = 45979.713148148148
$ws.Cells.Item(45, 4).Value = 45980.313275462962
$ws.Cells.Item(46, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(46, 2).Value = "904号直流"
$ws.Cells.Item(46, 3).Value = 45979.715844907405
$ws.Cells.Item(46, 4).Value = 45980.313275462962
$ws.Cells.Item(47, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(47, 2).Value = "505号直流"
$ws.Cells.Item(47, 3).Value = 45979.717314814814
$ws.Cells.Item(47, 4).Value = 45980.313275462962
$ws.Cells.Item(48, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(48, 2).Value = "307号直流"
$ws.Cells.Item(48, 3).Value = 45979.730312500003
$ws.Cells.Item(48, 4).Value = 45980.313275462962
$ws.Cells.Item(49, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(49, 2).Value = "302号直流"
$ws.Cells.Item(49, 3).Value = 45979.749525462961
$ws.Cells.Item(49, 4).Value = 45980.313275462962
$ws.Cells.Item(50, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(50, 2).Value = "105号直流"
$ws.Cells.Item(50, 3).Value = 45979.768703703703
$ws.Cells.Item(50, 4).Value = 45980.313275462962
$ws.Cells.Item(51, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(51, 2).Value = "205号直流"
$ws.Cells.Item(51, 3).Value = 45979.773217592592
$ws.Cells.Item(51, 4).Value = 45980.313275462962
$ws.Cells.Item(52, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(52, 2).Value = "804号直流"
$ws.Cells.Item(52, 3).Value = 45979.789699074077
$ws.Cells.Item(52, 4).Value = 45980.313275462962

# Update the saved selection to match the author's cursor position (E9)
$ws.Range("E9").Select()
